# outputs-HGR-r202/test-g__CAG-988_split_pruned.xlsx
# Re-exported prediction scores: column B ("1-s__CAG-988 sp003149915")
# on the "quadratic-svm-score" sheet moves from placeholder 1s to the
# actual computed decimal scores for rows 2-7 (Row values stay intact).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quadratic-svm-score")
if (-not $ws) { $ws = $wb.ActiveSheet }

$ws.Range("B2").Value = 2.8009227248630744
$ws.Range("B3").Value = 1.5703667084334576
$ws.Range("B4").Value = 1.5756503015243268
$ws.Range("B5").Value = 3.353323903427377
$ws.Range("B6").Value = 1.6312028632968634
$ws.Range("B7").Value = 2.8386188475372798
